# amlak_info.xlsx - "Sort - add columns - add filters - improvement"
#
# The sheet gains a new column B ("متولی" / custodian), the old "مساخت"
# header has its typo fixed to "مساحت" (area), and the old three columns
# "مالک" / "Latitude" / "Longitude" are replaced by two new columns
# "نوع مالکت" (ownership type) and "مختصات" (coordinates). Net column
# count stays the same (A:N) because +1 (new متولی column) and -1
# (3 old columns collapsed into 2 new ones) cancel out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at B - this shifts the old B:N headers (محل..نوع کاربری)
#    one slot to the right (C:O) and carries their column widths with them,
#    matching the width shift seen for columns C..N in the new layout.
$ws.Columns.Item(2).Insert()

# Give the freshly inserted column B the same width as column A (8.5703125
# in the saved OOXML). ColumnWidth is pixel-quantised by Excel, so 7.65
# characters is the closest achievable input that rounds to that width.
$ws.Columns.Item(2).ColumnWidth = 7.65

# 2) After the insert, the old "مالک" / "Latitude" / "Longitude" trio now
#    sits at K / L / M. Delete one of them (the old "Longitude" column, M)
#    so that three old columns collapse down to the two new ones, keeping
#    the total column count at 14 (A:N) and sliding "کد کاربری" /
#    "نوع کاربری" back into M / N with their original widths intact.
$ws.Columns.Item(13).Delete()

# 3) Write the final header row text (this also fixes the "مساخت" ->
#    "مساحت" typo and relabels the remaining cells).
$ws.Range("A1").Value = "شناسه "
$ws.Range("B1").Value = "متولی"
$ws.Range("C1").Value = "محل"
$ws.Range("D1").Value = "اطلاعات تکیل شده؟"
$ws.Range("E1").Value = "مساحت"
$ws.Range("F1").Value = "نوع"
$ws.Range("G1").Value = "نام"
$ws.Range("H1").Value = "آدرس"
$ws.Range("I1").Value = "وضعیت فعلی"
$ws.Range("J1").Value = "ساختار"
$ws.Range("K1").Value = "نوع مالکت"
$ws.Range("L1").Value = "مختصات"
$ws.Range("M1").Value = "کد کاربری"
$ws.Range("N1").Value = "نوع کاربری"

# 4) Move the saved selection from M3 to M4, as in the new workbook.
$ws.Range("M4").Select()
